$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 62, shifting existing rows 62-142 down to 63-143
$ws.Rows.Item(62).Insert()

# Populate the newly inserted row 62 with the new data record
$ws.Cells.Item(62, 1).Value = 11
$ws.Cells.Item(62, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(62, 3).Value = "Bíobío"
$ws.Cells.Item(62, 4).Value = 44799
$ws.Cells.Item(62, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(62, 5).Value = 8
$ws.Cells.Item(62, 6).Value = "Fruta"
$ws.Cells.Item(62, 7).Value = 100102
$ws.Cells.Item(62, 8).Value = "Cítricos"
$ws.Cells.Item(62, 9).Value = 100102004
$ws.Cells.Item(62, 10).Value = "Mandarina"
$ws.Cells.Item(62, 11).Value = "Murcott"
$ws.Cells.Item(62, 12).Value = "Primera"
$ws.Cells.Item(62, 13).Value = 220
$ws.Cells.Item(62, 14).Value = 7000
$ws.Cells.Item(62, 15).Value = 7500
$ws.Cells.Item(62, 16).Value = 7273
$ws.Cells.Item(62, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(62, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(62, 19).Value = 727
$ws.Cells.Item(62, 20).Value = 10
